# Weekly update: a new price record for "Achicoria" (Vega Modelo de Temuco)
# is inserted as the new row 44, pushing the existing rows 44-126 down to
# 45-127 (dimension grows from A1:R126 to A1:R127).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 44; this shifts rows 44..126
# down to 45..127 and carries their formatting (e.g. the date style on
# column D) along with them.
$ws.Rows("44:44").Insert()

# Fill in the new record in row 44.
$ws.Cells.Item(44, 1).Value = 10
$ws.Cells.Item(44, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(44, 3).Value = "La Araucanía"
$ws.Cells.Item(44, 4).Value = 45125
$ws.Cells.Item(44, 5).Value = 9
$ws.Cells.Item(44, 6).Value = 100112010
$ws.Cells.Item(44, 7).Value = "Achicoria"
$ws.Cells.Item(44, 8).Value = "Sin especificar"
$ws.Cells.Item(44, 9).Value = "Primera"
$ws.Cells.Item(44, 10).Value = 100
$ws.Cells.Item(44, 11).Value = 9000
$ws.Cells.Item(44, 12).Value = 9000
$ws.Cells.Item(44, 13).Value = 9000
$ws.Cells.Item(44, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(44, 15).Value = "Región Metropolitana"
$ws.Cells.Item(44, 16).Value = 500
$ws.Cells.Item(44, 17).Value = 18
$ws.Cells.Item(44, 18).Value = "Hortaliza"
